# Apply the v1.2.3 -> v1.2.4 changes to the UC007 test suite workbook.
#
# 1) Extend the "SYSTEM Recupera..." expected-result text (step 1 of every
#    test case, rows 10/19/27/35/43 in column D) with an extra sentence.
# 2) Rotate the "second step" content of test cases TC2/TC3/TC4 so that:
#      new TC2 second step = old TC4 second step (assign/unassign AP)
#      new TC3 second step = old TC2 second step (filter by user)
#      new TC4 second step = old TC3 second step (authorize payment)
#
# Note: reading via the `.Value` property of a Range in this runtime can
# return a bogus descriptor string, so `.Value2` is used for all reads
# (and, for consistency, all writes too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Append extra sentence to the recurring "SYSTEM Recupera..." text ---
$oldRecupera = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo numero de diarias em ordem crescente."
$newRecupera = $oldRecupera + " Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."

$recuperaRows = @(10, 19, 27, 35, 43)
foreach ($r in $recuperaRows) {
    $cell = $ws.Cells.Item($r, 4)  # column D
    if ($cell.Value2 -eq $oldRecupera) {
        $cell.Value2 = $newRecupera
    }
}

# --- 2) Rotate the second-step (row offset) content among TC2, TC3, TC4 ---

# Old content, captured before any mutation.
$filterStep   = $ws.Range("B20").Value2
$filterResult = $ws.Range("D20").Value2

$authorizeStep   = $ws.Range("B28").Value2
$authorizeResult = $ws.Range("D28").Value2

$assignStep   = $ws.Range("B36").Value2
$assignResult = $ws.Range("D36").Value2

# New TC2 second step <- old TC4 second step (assign/unassign AP)
$ws.Range("B20").Value2 = $assignStep
$ws.Range("D20").Value2 = $assignResult

# New TC3 second step <- old TC2 second step (filter by user)
$ws.Range("B28").Value2 = $filterStep
$ws.Range("D28").Value2 = $filterResult

# New TC4 second step <- old TC3 second step (authorize payment)
$ws.Range("B36").Value2 = $authorizeStep
$ws.Range("D36").Value2 = $authorizeResult
